$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04339299999999999
$ws.Range("H2").Value = 0.130179
$ws.Range("I2").Value = 0.0698021577815419
$ws.Range("J2").Value = 0.0698021577815419
$ws.Range("M2").Value = 0.7285076666666667
$ws.Range("N2").Value = 2.185523
$ws.Range("O2").Value = 0.1384760821597099
$ws.Range("P2").Value = 0.1384760821597099
$ws.Range("Q2").Value = 0.03161213317966666
$ws.Range("R2").Value = 0.284509198617
$ws.Range("S2").Value = 0.009665929335881833
$ws.Range("T2").Value = 0.009665929335881833
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04339299999999999
$ws.Range("H3").Value = 0.130179
$ws.Range("I3").Value = 0.0698021577815419
$ws.Range("J3").Value = 0.0698021577815419
$ws.Range("O3").Value = 0.6813230330092965
$ws.Range("P3").Value = 0.6813230330092966
$ws.Range("Q3").Value = 0.1555364227666666
$ws.Range("R3").Value = 1.3998278049
$ws.Range("S3").Value = 0.0475578178503136
$ws.Range("T3").Value = 0.04755781785031361
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04339299999999999
$ws.Range("H4").Value = 0.130179
$ws.Range("I4").Value = 0.0698021577815419
$ws.Range("J4").Value = 0.0698021577815419
$ws.Range("O4").Value = 0.1802008848309935
$ws.Range("P4").Value = 0.1802008848309935
$ws.Range("Q4").Value = 0.04113731614533333
$ws.Range("R4").Value = 0.370235845308
$ws.Range("S4").Value = 0.01257841059534647
$ws.Range("T4").Value = 0.01257841059534647
$ws.Range("I5").Value = 0.6764796878879081
$ws.Range("J5").Value = 0.6764796878879081
$ws.Range("M5").Value = 0.7285076666666667
$ws.Range("N5").Value = 2.185523
$ws.Range("O5").Value = 0.1384760821597099
$ws.Range("P5").Value = 0.1384760821597099
$ws.Range("Q5").Value = 0.3063653999605555
$ws.Range("R5").Value = 2.757288599645
$ws.Range("S5").Value = 0.0936762568393409
$ws.Range("T5").Value = 0.0936762568393409
$ws.Range("I6").Value = 0.6764796878879081
$ws.Range("J6").Value = 0.6764796878879081
$ws.Range("O6").Value = 0.6813230330092965
$ws.Range("P6").Value = 0.6813230330092966
$ws.Range("S6").Value = 0.4609011927209719
$ws.Range("T6").Value = 0.4609011927209719
$ws.Range("I7").Value = 0.6764796878879081
$ws.Range("J7").Value = 0.6764796878879081
$ws.Range("O7").Value = 0.1802008848309935
$ws.Range("P7").Value = 0.1802008848309935
$ws.Range("S7").Value = 0.1219022383275953
$ws.Range("T7").Value = 0.1219022383275953
$ws.Range("I8").Value = 0.2537181543305499
$ws.Range("J8").Value = 0.2537181543305499
$ws.Range("M8").Value = 0.7285076666666667
$ws.Range("N8").Value = 2.185523
$ws.Range("O8").Value = 0.1384760821597099
$ws.Range("P8").Value = 0.1384760821597099
$ws.Range("Q8").Value = 0.1149043573967778
$ws.Range("R8").Value = 1.034139216571
$ws.Range("S8").Value = 0.0351338959844872
$ws.Range("T8").Value = 0.0351338959844872
$ws.Range("I9").Value = 0.2537181543305499
$ws.Range("J9").Value = 0.2537181543305499
$ws.Range("O9").Value = 0.6813230330092965
$ws.Range("P9").Value = 0.6813230330092966
$ws.Range("S9").Value = 0.1728640224380111
$ws.Range("T9").Value = 0.1728640224380111
$ws.Range("I10").Value = 0.2537181543305499
$ws.Range("J10").Value = 0.2537181543305499
$ws.Range("O10").Value = 0.1802008848309935
$ws.Range("P10").Value = 0.1802008848309935
$ws.Range("S10").Value = 0.04572023590805166
$ws.Range("T10").Value = 0.04572023590805166
